$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D width (target stored width 24.17; Excel quantizes to the
# nearest 1/6 character width, so 23.3 is the ColumnWidth that rounds to it) ---
$ws.Columns.Item(4).ColumnWidth = 23.3

# --- Row 3: always_not_taken ---
$ws.Range("D3").Value = '75%(for real random)'
$ws.Range("E3").Value = 'Main loop is always taken (unconditional jump) => 50% of all prediction is mistaken. Tested branch is branch that uses current bit to compare them with zero. => if it was real random we would have 50%(of ½ of all predictions) mistakes(only on this branch) => 50+25 = 75%. But our pattern is 01111100010110100001100001001110 (15 – true and 17 false) => 15/32 * ½ + 50% = 23.4 + 50 = 73.4'

# --- Row 4: always_taken ---
$ws.Range("B4").Value = "always_taken"
$ws.Range("D4").Value = '25%(for real random)'
$ws.Range("E4").Value = 'Main loop is always taken (unconditional jump) => 50% of all prediction is mistaken. Tested branch is branch that uses current bit to compare them with zero. => if it was real random we would have 50%(of ½ of all predictions) mistakes(only on this branch) => ½ * 1/2 = 25%. But our pattern is 01111100010110100001100001001110 (15 – true and 17 false) => 15/32 * ½ = 26.6%'

# --- Row 5: backward_jumps ---
$ws.Range("B5").Value = "backward_jumps"
$ws.Range("D5").Value = '25%(for real random)'
$ws.Range("E5").Value = 'Main loop is unconditional backward jump => 1/2 of all predictions is predicted. Tested branch is branch(forward jump and this predictor doesn’t work) that uses current bit to compare them with zero. => if it was real random we would have 50%(of ½ of all predictions) mistakes(only on this branch) => ½ * 1/2 = 25%. But our pattern is 01111100010110100001100001001110 (15 – true and 17 false) => 15/32 * ½ = 26.6%'
$ws.Range("E5").Style = "Normal"

# --- Row 6: saturating_one_bit ---
$ws.Range("B6").Value = "saturating_one_bit"
$ws.Range("D6").Value = '25%(for real random)'
$ws.Range("E6").Value = 'Main loop is always taken and predictor knows it after the first time => 1/2 predictions is true. Tested branch is branch that uses current bit to compare them with zero. => if it was real random we would have 50%(of ½ of all predictions) mistakes(only on this branch) => ½ * 1/2 = 25%. But we have 01111100010110100001100001001110 with 18 right predictions and 14 misses => 14/32 = 21.9%'
$ws.Range("E6").Style = "Normal"

# --- Row 7: saturating_two_bits ---
$ws.Range("B7").Value = "saturating_two_bits"
$ws.Range("D7").Value = '25%(for real random)'
$ws.Range("E7").Value = 'Main loop is always taken and predictor knows it after the first time => 1/2 predictions is true. Tested branch is branch that uses current bit to compare them with zero. => if it was real random we would have 50%(of ½ of all predictions) mistakes(only on this branch) => ½ * 1/2 = 25%. But we have 01111100010110100001100001001110 with 14 right predictions and 18 misses => 18/32 * ½  = 28.1%(one bit is better because it is simpler to make it “hot”, SN → WEAKLY T is 2 stages, but N → T is only one stage) (in other situation with other random pattern two-bits is better, but it’s problems of our method of producing random numbers)'

# --- Row 8: adaptive_two_levels ---
$ws.Range("B8").Value = "adaptive_two_levels"
$ws.Range("D8").Value = '25%(for real random)'
$ws.Range("E8").Value = 'Main loop is always taken and predictor knows it after the first time => 1/2 predictions is true. Tested branch is branch that uses current bit to compare them with zero. => if it was real random we would have 50%(of ½ of all predictions) mistakes(only on this branch) => ½ * 1/2 = 25%. But we have 01111100010110100001100001001110 with 14 right predictions and 18 misses => 18/32 * ½  = 28.1%(it is worse than some simpler predictor because two-levels use 2bit pattern that can remember bad  predictions. But this effect can be seen only on special random patterns (in other situation with other random pattern two-levels is better, but it’s problems of our method of producing random numbers)'

# --- Sheet view ---
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$ws.Range("E20").Select() | Out-Null
